# Auto-generated: apply scheduled-runner market-data refresh to Excalibur_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 218.35715
$ws.Range("I9").Value = 150.58333
$ws.Range("K9").Value = 150.58333
$ws.Range("M9").Value = 18.41667000000001
# Row 32
$ws.Range("H32").Value = 3449
$ws.Range("I32").Value = 4499
$ws.Range("J32").Value = 2399
$ws.Range("K32").Value = 4499
$ws.Range("L32").Value = 2399
$ws.Range("M32").Value = -4173
$ws.Range("N32").Value = -3051
# Row 34
$ws.Range("H34").Value = 9138.571
$ws.Range("I34").Value = 9138.571
$ws.Range("K34").Value = 9138.571
$ws.Range("M34").Value = -8935.571
# Row 36
$ws.Range("H36").Value = 9138.571
$ws.Range("I36").Value = 9138.571
$ws.Range("K36").Value = 9138.571
$ws.Range("M36").Value = -8423.571
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
# Row 86
$ws.Range("H86").Value = 1901.7142
$ws.Range("I86").Value = 1115.8572
$ws.Range("J86").Value = 2687.5715
$ws.Range("K86").Value = 1115.8572
$ws.Range("L86").Value = 2687.5715
$ws.Range("M86").Value = 7.142800000000079
$ws.Range("N86").Value = -4933.5715
# Row 89
$ws.Range("H89").Value = 1901.7142
$ws.Range("I89").Value = 1115.8572
$ws.Range("J89").Value = 2687.5715
$ws.Range("K89").Value = 5579.286
$ws.Range("L89").Value = 13437.8575
$ws.Range("M89").Value = 36.71399999999994
$ws.Range("N89").Value = -24669.8575
# Row 106
$ws.Range("H106").Value = 2703.6667
$ws.Range("I106").Value = 2651.6365
$ws.Range("K106").Value = 2651.6365
$ws.Range("M106").Value = -2020.6365

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 68673.266
$ws.Range("I37").Value = 43750
$ws.Range("J37").Value = 77736.27
$ws.Range("K37").Value = 43750
$ws.Range("L37").Value = 77736.27
$ws.Range("M37").Value = -43477
$ws.Range("N37").Value = -78282.27
# Row 61
$ws.Range("H61").Value = 8334111
$ws.Range("I61").Value = 8334111
$ws.Range("K61").Value = 8334111
$ws.Range("M61").Value = -8333899
# Row 74
$ws.Range("H74").Value = 4489.9473
$ws.Range("I74").Value = 1370.2222
$ws.Range("K74").Value = 1370.2222
$ws.Range("M74").Value = -496.2221999999999
# Row 77
$ws.Range("H77").Value = 4489.9473
$ws.Range("I77").Value = 1370.2222
$ws.Range("K77").Value = 6851.111
$ws.Range("M77").Value = -2483.111
# Row 110
$ws.Range("H110").Value = 2616.8845
$ws.Range("I110").Value = 2378.0952
$ws.Range("J110").Value = 3619.8
$ws.Range("K110").Value = 2378.0952
$ws.Range("L110").Value = 3619.8
$ws.Range("M110").Value = -333.0952000000002
$ws.Range("N110").Value = -7709.8
# Row 136
$ws.Range("H136").Value = 8334111
$ws.Range("I136").Value = 8334111
$ws.Range("K136").Value = 25002333
$ws.Range("M136").Value = -24999783

$ws = $wb.Worksheets.Item("BSM")
# Row 48
$ws.Range("H48").Value = 400592
$ws.Range("J48").Value = 400592
$ws.Range("L48").Value = 400592
$ws.Range("N48").Value = -401422
# Row 105
$ws.Range("H105").Value = 8146.143
$ws.Range("I105").Value = 8146.143
$ws.Range("K105").Value = 8146.143
$ws.Range("M105").Value = -6399.143

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 679.2857
$ws.Range("I22").Value = 711
$ws.Range("K22").Value = 711
$ws.Range("M22").Value = -361
# Row 58
$ws.Range("H58").Value = 1553422.2
$ws.Range("J58").Value = 13001.8
$ws.Range("L58").Value = 13001.8
$ws.Range("N58").Value = -13407.8
# Row 62
$ws.Range("H62").Value = 3401.8
$ws.Range("I62").Value = 2799
$ws.Range("J62").Value = 3552.5
$ws.Range("K62").Value = 2799
$ws.Range("L62").Value = 3552.5
$ws.Range("M62").Value = -2175
$ws.Range("N62").Value = -4800.5
# Row 65
$ws.Range("H65").Value = 3401.8
$ws.Range("I65").Value = 2799
$ws.Range("J65").Value = 3552.5
$ws.Range("K65").Value = 13995
$ws.Range("L65").Value = 17762.5
$ws.Range("M65").Value = -10875
$ws.Range("N65").Value = -24002.5
# Row 136
$ws.Range("H136").Value = 1553422.2
$ws.Range("J136").Value = 13001.8
$ws.Range("L136").Value = 39005.39999999999
$ws.Range("N136").Value = -44105.39999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 115
$ws.Range("H115").Value = 2167.3333
$ws.Range("I115").Value = 2167.3333
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 6501.999899999999
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -5326.999899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 323976.7
$ws.Range("I80").Value = 515512.1
$ws.Range("J80").Value = 4751
$ws.Range("K80").Value = 515512.1
$ws.Range("L80").Value = 4751
$ws.Range("M80").Value = -514514.1
$ws.Range("N80").Value = -6747
# Row 83
$ws.Range("H83").Value = 323976.7
$ws.Range("I83").Value = 515512.1
$ws.Range("J83").Value = 4751
$ws.Range("K83").Value = 2577560.5
$ws.Range("L83").Value = 23755
$ws.Range("M83").Value = -2572568.5
$ws.Range("N83").Value = -33739
# Row 126
$ws.Range("H126").Value = 759713.7
$ws.Range("I126").Value = 1112744.9
$ws.Range("J126").Value = 3218.1428
$ws.Range("K126").Value = 3338234.7
$ws.Range("L126").Value = 9654.428400000001
$ws.Range("M126").Value = -3335764.7
$ws.Range("N126").Value = -14594.4284

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3309.7273
$ws.Range("I40").Value = 3256.9062
$ws.Range("K40").Value = 3256.9062
$ws.Range("M40").Value = -3120.9062
# Row 42
$ws.Range("H42").Value = 17674.666
$ws.Range("I42").Value = 17674.666
$ws.Range("K42").Value = 17674.666
$ws.Range("M42").Value = -17111.666
# Row 43
$ws.Range("H43").Value = 343666.66
$ws.Range("I43").Value = 6000
$ws.Range("J43").Value = 512500
$ws.Range("K43").Value = 6000
$ws.Range("L43").Value = 512500
$ws.Range("M43").Value = -5807
$ws.Range("N43").Value = -512886
# Row 46
$ws.Range("H46").Value = 6964.5835
$ws.Range("I46").Value = 9547.25
$ws.Range("K46").Value = 9547.25
$ws.Range("M46").Value = -9359.25
# Row 49
$ws.Range("H49").Value = 17674.666
$ws.Range("I49").Value = 17674.666
$ws.Range("K49").Value = 17674.666
$ws.Range("M49").Value = -17527.666
# Row 82
$ws.Range("H82").Value = 1013.4286
$ws.Range("I82").Value = 898.75
$ws.Range("K82").Value = 898.75
$ws.Range("M82").Value = -537.75
# Row 85
$ws.Range("H85").Value = 1013.4286
$ws.Range("I85").Value = 898.75
$ws.Range("K85").Value = 898.75
$ws.Range("M85").Value = 349.25
# Row 93
$ws.Range("H93").Value = 2949.6667
$ws.Range("I93").Value = 2949.6667
$ws.Range("K93").Value = 2949.6667
$ws.Range("M93").Value = -1701.6667
# Row 132
$ws.Range("H132").Value = 3166831
$ws.Range("I132").Value = 3869015.8
$ws.Range("K132").Value = 11607047.4
$ws.Range("M132").Value = -11604517.4

$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = 0
# Row 107
$ws.Range("H107").Value = 1250.125
$ws.Range("I107").Value = 769.619
$ws.Range("J107").Value = 2167.4546
$ws.Range("K107").Value = 2308.857
$ws.Range("L107").Value = 6502.3638
$ws.Range("M107").Value = -388.857
$ws.Range("N107").Value = -10342.3638
# Row 122
$ws.Range("H122").Value = 2757.9375
$ws.Range("I122").Value = 2383.8518
$ws.Range("K122").Value = 7151.555399999999
$ws.Range("M122").Value = -4701.555399999999
# Row 126
$ws.Range("H126").Value = 4658.0435
$ws.Range("I126").Value = 4356.75
$ws.Range("K126").Value = 13070.25
$ws.Range("M126").Value = -10600.25
# Row 132
$ws.Range("H132").Value = 8752879
$ws.Range("I132").Value = 11182817
$ws.Range("J132").Value = 5100
$ws.Range("K132").Value = 33548451
$ws.Range("L132").Value = 15300
$ws.Range("M132").Value = -33545921
$ws.Range("N132").Value = -20360
# Row 136
$ws.Range("H136").Value = 6586.7617
$ws.Range("I136").Value = 6809.222
$ws.Range("K136").Value = 20427.666
$ws.Range("M136").Value = -17877.666
